{"js": "// Bug-fix update: rename the \"Informacao espec\u00edfica...\" bug item to\n// \"Numero acima de \"4.20e+3\"\" and drop the now-duplicate bullet that used\n// to hold that text, then keep Word's \"_GoBack\" bookmark anchored right\n// after the \"Desabilitar anos...\" bullet (i.e. on the paragraph that now\n// follows it), just like the saved document does.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the relevant paragraphs by their current text so the script is\n// resilient to exact indices.\nlet infoIdx = null;\nlet numeroIdx = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (infoIdx === null && t.indexOf(\"Informacao espec\u00edfica no hover das divis\u00f5es\") !== -1) {\n    infoIdx = i;\n  } else if (numeroIdx === null && t.indexOf(\"Numero acima de\") !== -1) {\n    numeroIdx = i;\n  }\n}\n\n// 1) Replace the text of the \"Informacao espec\u00edfica...\" bullet.\nparagraphs.items[infoIdx].insertText(\"Numero acima de \\u201c4.20e+3\\u201d\", Word.InsertLocation.replace);\n\n// 2) Remove the paragraph that used to hold \"Numero acima de ...\" (now a\n// duplicate of the text we just wrote above).\nparagraphs.items[numeroIdx].delete();\n\nawait context.sync();\n\n// 3) Move the \"_GoBack\" bookmark off of the \"Desabilitar anos...\" bullet\n// and onto the paragraph right after it (matching Word's own bookkeeping\n// of the last-edit location).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst refreshedParagraphs = context.document.body.paragraphs;\nrefreshedParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet afterDesabilitarIdx = null;\nfor (let i = 0; i < refreshedParagraphs.items.length; i++) {\n  const t = refreshedParagraphs.items[i].text;\n  if (t.indexOf(\"Desabilitar anos que n\u00e3o geram dados\") !== -1) {\n    afterDesabilitarIdx = i + 1;\n    break;\n  }\n}\n\nconst bookmarkRange = refreshedParagraphs.items[afterDesabilitarIdx].getRange(Word.RangeLocation.start);\nbookmarkRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Bug-fix update: rename the \"Informacao espec\u00edfica...\" bug item to\n# \"Numero acima de \"4.20e+3\"\" and drop the now-duplicate bullet that used\n# to hold that text, then keep Word's \"_GoBack\" bookmark anchored right\n# after the \"Desabilitar anos...\" bullet (i.e. on the paragraph that now\n# follows it), just like the saved document does.\n\n$d = $word.ActiveDocument\n\n$quoteOpen = [char]0x201C\n$quoteClose = [char]0x201D\n\n# Locate the relevant paragraphs by their current text so the script is\n# resilient to exact indices.\n$infoParaIndex = $null\n$numeroParaIndex = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($infoParaIndex -eq $null -and $t -like \"*Informacao espec\u00edfica no hover das divis\u00f5es*\") {\n        $infoParaIndex = $i\n    } elseif ($numeroParaIndex -eq $null -and $t -like \"*Numero acima de*\") {\n        $numeroParaIndex = $i\n    }\n}\n\n# 1) Replace the text of the \"Informacao espec\u00edfica...\" bullet.\n$d.Paragraphs.Item($infoParaIndex).Range.Text = \"Numero acima de \" + $quoteOpen + \"4.20e+3\" + $quoteClose\n\n# 2) Remove the paragraph that used to hold \"Numero acima de ...\" (now a\n# duplicate of the text we just wrote above).\n$d.Paragraphs.Item($numeroParaIndex).Range.Delete()\n\n# 3) Move the \"_GoBack\" bookmark off of the \"Desabilitar anos...\" bullet\n# and onto the paragraph right after it (matching Word's own bookkeeping\n# of the last-edit location). Re-locate \"Desabilitar anos...\" after the\n# deletion above since paragraph indices shifted by one.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$desabilitarParaIndex = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($desabilitarParaIndex -eq $null -and $t -like \"*Desabilitar anos que n\u00e3o geram dados*\") {\n        $desabilitarParaIndex = $i\n        break\n    }\n}\n\n$afterDesabilitarIndex = $desabilitarParaIndex + 1\n$targetRange = $d.Paragraphs.Item($afterDesabilitarIndex).Range\n$targetRange.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $targetRange)\n"}
